$wb = $excel.ActiveWorkbook

# --- Insert the new "report_locations" sheet between "input" and "who_dummies" ---
$inputSheet = $wb.Worksheets.Item("input")
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $inputSheet)
$ws.Name = "report_locations"

# Header row: report / dir (column B filled first, matches author's original entry order)
$ws.Cells.Item(1, 1).Value = "report"
$ws.Cells.Item(1, 2).Value = "dir"

# Directory column (B2:B6) for each report
$paths = @(
    "C:/Users/hswerdfe/Documents/reports/qry",
    "C:/Users/hswerdfe/Documents/reports/DS",
    "C:/Users/hswerdfe/Documents/reports/HC",
    "C:/Users/hswerdfe/Documents/reports/STAT",
    "C:/Users/hswerdfe/Documents/reports/WHO"
)
for ($r = 0; $r -lt $paths.Length; $r++) {
    $ws.Cells.Item($r + 2, 2).Value = $paths[$r]
}

# Day-of-week "X" flags (columns C:I, rows 2:6) - every report runs every day
for ($r = 2; $r -le 6; $r++) {
    for ($c = 3; $c -le 9; $c++) {
        $ws.Cells.Item($r, $c).Value = "X"
    }
}

# Day-of-week header labels (row 1, columns C:I)
$days = @("Sun", "Mon", "Tue", "Wed", "Thu", "Fri", "Sat")
for ($c = 0; $c -lt $days.Length; $c++) {
    $ws.Cells.Item(1, $c + 3).Value = $days[$c]
}

# Report-name column (A2:A6)
$reports = @("qry_allcases", "Domestic surveillance ", "HCDaily", "STATCAN", "WHO")
for ($r = 0; $r -lt $reports.Length; $r++) {
    $ws.Cells.Item($r + 2, 1).Value = $reports[$r]
}

# Row 5 (STATCAN) has a slightly custom height in the target workbook
$ws.Rows.Item(5).RowHeight = 14.25

# --- View-state changes ---
# "input" sheet scrolls down and selects B284:B291, and is no longer the active tab
[void]$inputSheet.Activate()
[void]$inputSheet.Range("B284:B291").Select()

# "report_locations" becomes the active/selected sheet & cell
[void]$ws.Activate()
[void]$ws.Range("B2").Select()
